$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain stored as text
# (matching the workbook's original inlineStr/text cell type). Setting
# NumberFormat to "@" (Text) on each target cell before assigning its value
# prevents Excel from auto-converting the numeric-looking string to a number.

# --- Column D (Price) updates ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.27'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.68'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.421'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05767'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.434'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.331'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8113'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8991'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1446'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07343'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03143'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.02991'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09415'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001578'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04820'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0005851'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.006139'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.004062'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0009956'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0001500'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.748'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.199'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.3279'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1330'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.178'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003159'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03917'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006779'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1073'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003200'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007322'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005643'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.3801'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1680'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'

# --- Other column updates (Coin name / Link / Volume(1h) label) ---
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E16').Value = '15CoinExTokenCET'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E17').Value = '16OneONE'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('E18').Value = '17TigerCashTCH'
$ws.Range('B19').Value = 'HotbitToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E19').Value = '18HotbitTokenHTB'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'NitroEx'
$ws.Range('C21').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E21').Value = '20NitroExNTX'
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E22').Value = '21LEOLEO'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'BitpandaEcosystemToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('E24').Value = '23BitpandaEcosystemTokenBEST'
$ws.Range('B25').Value = 'ProBitToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('E25').Value = '24ProBitTokenPROB'
$ws.Range('B26').Value = 'MCDex'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E26').Value = '25MCDexMCB'
$ws.Range('E43').Value = '42CEJICEJIBestin24h'
